$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3's data (name, place, number) into a new row 7,
# so the output can contain repeated names (two rows for "لؤي").
$ws.Range("A7").Value = $ws.Range("A3").Value2
$ws.Range("B7").Value = $ws.Range("B3").Value2
$ws.Range("C7").Value = $ws.Range("C3").Value2

# Copy row 3's cell formatting onto row 7 (reuse the existing style
# instead of fabricating a new one).
$ws.Range("A3:C3").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Update the active selection to reflect the newly edited row.
$ws.Range("B3").Select()
